$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Header updates
$ws1.Cells.Item(2,1).Value = "Última actualización: 15:44:42"
$ws1.Cells.Item(3,1).Value = "Total filas: 225"

# Swap the A/C/D values between adjacent tie rows (re-sort artifact of the scrape)
$ws1.Cells.Item(62,1).Value = "07:49:32"
$ws1.Cells.Item(62,3).Value = "14_ABASTO"
$ws1.Cells.Item(62,4).Value = 88
$ws1.Cells.Item(63,1).Value = "08:38:24"
$ws1.Cells.Item(63,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(63,4).Value = 39

$ws1.Cells.Item(118,1).Value = "10:36:50"
$ws1.Cells.Item(118,3).Value = "225_GOMEZ"
$ws1.Cells.Item(118,4).Value = 76
$ws1.Cells.Item(119,1).Value = "11:33:52"
$ws1.Cells.Item(119,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(119,4).Value = 19

$ws1.Cells.Item(133,1).Value = "11:33:52"
$ws1.Cells.Item(133,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(133,4).Value = 59
$ws1.Cells.Item(134,1).Value = "10:36:50"
$ws1.Cells.Item(134,3).Value = "14_ABASTO"
$ws1.Cells.Item(134,4).Value = 116

$ws1.Cells.Item(147,1).Value = "11:33:52"
$ws1.Cells.Item(147,3).Value = "215C_EL PATO"
$ws1.Cells.Item(147,4).Value = 90
$ws1.Cells.Item(148,1).Value = "11:13:15"
$ws1.Cells.Item(148,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(148,4).Value = 110

$ws1.Cells.Item(160,1).Value = "12:33:02"
$ws1.Cells.Item(160,3).Value = "14_ABASTO"
$ws1.Cells.Item(160,4).Value = 60
$ws1.Cells.Item(161,1).Value = "11:46:32"
$ws1.Cells.Item(161,3).Value = "215A_EL PATO"
$ws1.Cells.Item(161,4).Value = 107

# New scrape (15:44:42) adds 6 rows, merged in Hora_Llegada order.
# Insert a row before old row 217 (Hora_Llegada 16:30) for the new 16:29 arrival.
$ws1.Rows.Item(217).Insert()
$ws1.Cells.Item(217,1).Value = "15:44:42"
$ws1.Cells.Item(217,2).Value = "16:29"
$ws1.Cells.Item(217,3).Value = "14_ABASTO"
$ws1.Cells.Item(217,4).Value = 45
$ws1.Cells.Item(217,5).Value = "LP1912"

# Insert a row before (now-shifted) old row 223 (Hora_Llegada 17:07) for the new 17:02 arrival.
$ws1.Rows.Item(224).Insert()
$ws1.Cells.Item(224,1).Value = "15:44:42"
$ws1.Cells.Item(224,2).Value = "17:02"
$ws1.Cells.Item(224,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(224,4).Value = 78
$ws1.Cells.Item(224,5).Value = "LP1912"

# Append the remaining 4 new arrivals at the end of the sheet.
$ws1.Cells.Item(227,1).Value = "15:44:42"
$ws1.Cells.Item(227,2).Value = "17:21"
$ws1.Cells.Item(227,3).Value = "15X38_ABASTO"
$ws1.Cells.Item(227,4).Value = 97
$ws1.Cells.Item(227,5).Value = "LP1912"

$ws1.Cells.Item(228,1).Value = "15:44:42"
$ws1.Cells.Item(228,2).Value = "17:34"
$ws1.Cells.Item(228,3).Value = "17_ROMERO"
$ws1.Cells.Item(228,4).Value = 110
$ws1.Cells.Item(228,5).Value = "LP1912"

$ws1.Cells.Item(229,1).Value = "15:44:42"
$ws1.Cells.Item(229,2).Value = "17:36"
$ws1.Cells.Item(229,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(229,4).Value = 112
$ws1.Cells.Item(229,5).Value = "LP1912"

$ws1.Cells.Item(230,1).Value = "15:44:42"
$ws1.Cells.Item(230,2).Value = "17:38"
$ws1.Cells.Item(230,3).Value = "215B_EL PATO"
$ws1.Cells.Item(230,4).Value = 114
$ws1.Cells.Item(230,5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(2,1).Value = "Última actualización: 15:44:42"
$ws2.Cells.Item(3,1).Value = "Total filas: 38"

$ws2.Cells.Item(43,1).Value = "15:44:42"
$ws2.Cells.Item(43,2).Value = "17:38"
$ws2.Cells.Item(43,3).Value = "215B_EL PATO"
$ws2.Cells.Item(43,4).Value = 114
$ws2.Cells.Item(43,5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Item(2,1).Value = "Última actualización: 15:44:42"
